$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4

$ws.Range("C3").Value = 18
$ws.Range("D3").Value = 12
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1

$ws.Range("C4").Value = 63
$ws.Range("D4").Value = 45
$ws.Range("E4").Value = 6

$ws.Range("C5").Value = 22
$ws.Range("D5").Value = 21
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1

$ws.Range("C6").Value = 35
$ws.Range("D6").Value = 37
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0

$ws.Range("C7").Value = 54
$ws.Range("D7").Value = 40
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 2

$ws.Range("C8").Value = 50
$ws.Range("D8").Value = 41
$ws.Range("E8").Value = 5

$ws.Range("C9").Value = 25
$ws.Range("D9").Value = 17
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 0

$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0

$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 8
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0

$ws.Range("C12").Value = 32
$ws.Range("D12").Value = 23
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 1

$ws.Range("C13").Value = 33
$ws.Range("D13").Value = 34
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1

$ws.Range("C14").Value = 56
$ws.Range("D14").Value = 42
$ws.Range("E14").Value = 8

$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 0

$ws.Range("C16").Value = 74
$ws.Range("D16").Value = 45
$ws.Range("E16").Value = 12
$ws.Range("F16").Value = 1
